{"js": "// Update the worksheet date and the 25 three-digit-by-one-digit\n// multiplication prompts to a newly generated set of problems.\nconst replacements = [\n  [\"2025-10-17 Friday\", \"2025-10-18 Saturday\"],\n  [\"212\u00d78=\", \"478\u00d74=\"],\n  [\"424\u00d79=\", \"892\u00d77=\"],\n  [\"616\u00d73=\", \"633\u00d78=\"],\n  [\"926\u00d73=\", \"791\u00d77=\"],\n  [\"207\u00d74=\", \"389\u00d78=\"],\n  [\"582\u00d79=\", \"619\u00d77=\"],\n  [\"935\u00d78=\", \"222\u00d73=\"],\n  [\"817\u00d75=\", \"992\u00d74=\"],\n  [\"341\u00d72=\", \"613\u00d74=\"],\n  [\"967\u00d74=\", \"739\u00d75=\"],\n  [\"153\u00d72=\", \"190\u00d72=\"],\n  [\"318\u00d73=\", \"517\u00d76=\"],\n  [\"390\u00d74=\", \"752\u00d76=\"],\n  [\"418\u00d78=\", \"284\u00d72=\"],\n  [\"422\u00d78=\", \"488\u00d77=\"],\n  [\"429\u00d77=\", \"437\u00d78=\"],\n  [\"186\u00d72=\", \"788\u00d73=\"],\n  [\"786\u00d74=\", \"227\u00d75=\"],\n  [\"189\u00d76=\", \"564\u00d77=\"],\n  [\"333\u00d76=\", \"376\u00d73=\"],\n  [\"726\u00d78=\", \"114\u00d78=\"],\n  [\"940\u00d75=\", \"255\u00d78=\"],\n  [\"871\u00d75=\", \"407\u00d76=\"],\n  [\"919\u00d74=\", \"816\u00d76=\"],\n  [\"457\u00d72=\", \"841\u00d73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the 25 three-digit-by-one-digit\n# multiplication prompts to a newly generated set of problems.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-10-17 Friday\", \"2025-10-18 Saturday\"),\n    @(\"212\u00d78=\", \"478\u00d74=\"),\n    @(\"424\u00d79=\", \"892\u00d77=\"),\n    @(\"616\u00d73=\", \"633\u00d78=\"),\n    @(\"926\u00d73=\", \"791\u00d77=\"),\n    @(\"207\u00d74=\", \"389\u00d78=\"),\n    @(\"582\u00d79=\", \"619\u00d77=\"),\n    @(\"935\u00d78=\", \"222\u00d73=\"),\n    @(\"817\u00d75=\", \"992\u00d74=\"),\n    @(\"341\u00d72=\", \"613\u00d74=\"),\n    @(\"967\u00d74=\", \"739\u00d75=\"),\n    @(\"153\u00d72=\", \"190\u00d72=\"),\n    @(\"318\u00d73=\", \"517\u00d76=\"),\n    @(\"390\u00d74=\", \"752\u00d76=\"),\n    @(\"418\u00d78=\", \"284\u00d72=\"),\n    @(\"422\u00d78=\", \"488\u00d77=\"),\n    @(\"429\u00d77=\", \"437\u00d78=\"),\n    @(\"186\u00d72=\", \"788\u00d73=\"),\n    @(\"786\u00d74=\", \"227\u00d75=\"),\n    @(\"189\u00d76=\", \"564\u00d77=\"),\n    @(\"333\u00d76=\", \"376\u00d73=\"),\n    @(\"726\u00d78=\", \"114\u00d78=\"),\n    @(\"940\u00d75=\", \"255\u00d78=\"),\n    @(\"871\u00d75=\", \"407\u00d76=\"),\n    @(\"919\u00d74=\", \"816\u00d76=\"),\n    @(\"457\u00d72=\", \"841\u00d73=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
